$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing the existing rows 16:124 down to 17:125.
# Excel's Rows.Insert shifts all cell data (and the sheet dimension) down
# automatically, which matches the shift seen across the whole diff.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new week's data point.
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "Femacal de La Calera"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44462
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 100112010
$ws.Range("G16").Value = "Achicoria"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 5500
$ws.Range("M16").Value = 5269
$ws.Range("N16").Value = "$/caja 16 unidades"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 329
$ws.Range("Q16").Value = 16
$ws.Range("R16").Value = "Hortaliza"
